# This script updates the cryptocurrency price/volume table with
# refreshed figures produced by the "Updated symbol list" GitHub Actions
# job. Column D (Price) and Column E (Volume(1h)) are text-formatted
# values (e.g. "310.10", "0.58%"), so we prefix each assignment with a
# leading apostrophe to force Excel to keep them as text instead of
# auto-converting them to numeric/percentage values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = "310.10"; E = "0.58%" }
    @{ Row = 3; D = "39.41"; E = "1.53%" }
    @{ Row = 4; D = "5.116"; E = "0.21%" }
    @{ Row = 5; D = "0.08112"; E = "-0.28%" }
    @{ Row = 6; D = "1.952"; E = "-1.02%" }
    @{ Row = 7; D = "8.162"; E = "2.91%" }
    @{ Row = 8; D = "0.9291"; E = "0.11%" }
    @{ Row = 9; D = "0.1409"; E = "-2.23%" }
    @{ Row = 10; D = "0.1932"; E = "-1.09%" }
    @{ Row = 11; D = "0.09066"; E = "-0.78%" }
    @{ Row = 12; D = $null; E = "-0.03%" }
    @{ Row = 13; D = "0.09819"; E = "-0.22%" }
    @{ Row = 14; D = "0.001397"; E = "-1.28%" }
    @{ Row = 15; D = "0.006102"; E = "3.56%" }
    @{ Row = 16; D = "3.913"; E = "8.87%" }
    @{ Row = 17; D = "4.244"; E = "1.06%" }
    @{ Row = 18; D = "3.322"; E = "-4.19%" }
    @{ Row = 19; D = $null; E = "0.23%" }
    @{ Row = 20; D = "0.1313"; E = "-0.06%" }
    @{ Row = 21; D = "4.732"; E = "-1.46%" }
    @{ Row = 22; D = "0.2427"; E = "-0.36%" }
    @{ Row = 23; D = "0.04368"; E = "-2.04%" }
    @{ Row = 24; D = "0.001230"; E = "-0.98%" }
    @{ Row = 25; D = "0.004804"; E = "-0.68%" }
    @{ Row = 26; D = "0.0001301"; E = "-0.12%" }
    @{ Row = 27; D = "0.0004005"; E = "-9.96%" }
    @{ Row = 39; D = "0.02150"; E = "2.39%" }
    @{ Row = 40; D = "0.05090"; E = "-0.36%" }
    @{ Row = 41; D = "0.007426"; E = "-0.61%" }
    @{ Row = 42; D = "0.009765"; E = "-3.25%" }
    @{ Row = 43; D = "0.1363"; E = "-0.07%" }
    @{ Row = 44; D = "0.002132"; E = "-0.58%" }
    @{ Row = 45; D = "0.009047"; E = "-13.69%" }
    @{ Row = 46; D = "0.00006408"; E = "3.17%" }
    @{ Row = 47; D = $null; E = "-0.11%" }
    @{ Row = 48; D = "0.001001"; E = "-37.57%" }
    @{ Row = 49; D = "0.002576"; E = "-15.77%" }
    @{ Row = 50; D = "0.00002102"; E = "-0.11%" }
    @{ Row = 51; D = "0.0002002"; E = "-0.11%" }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $ws.Cells.Item($u.Row, 4).Value = "'" + $u.D
    }
    $ws.Cells.Item($u.Row, 5).Value = "'" + $u.E
}

